# CableChambers.docx edit:
#  1. Merge several runs of split sentences back into single runs.
#  2. Replace the three embedded <w:drawing> pictures with hyperlink runs
#     that display the image's source URL as a clickable link.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Text merges (collapse runs that were unnecessarily split)
# ---------------------------------------------------------------

$mergeRange = $d.Content

$m1 = "Cable chambers integrated within the main building are excluded from GFA if it does not exceed 2.5m in height"
$mergeRange.Find.Execute($m1, $true, $false, $false, $false, $false, $true, 1, $false, $m1, 2) | Out-Null

$m2 = [char]0x00A0 + "The structural beams shall also be included in the height computation."
$mergeRange.Find.Execute($m2, $true, $false, $false, $false, $false, $true, 1, $false, $m2, 2) | Out-Null

$m3 = "Standalone cable chambers that fulfil either criteria below are excluded from GFA:"
$mergeRange.Find.Execute($m3, $true, $false, $false, $false, $false, $true, 1, $false, $m3, 2) | Out-Null

$m4 = [char]0x00A0 + "of the cable chamber, measured from the floor to the underside of the upper floor shall not exceed 2.0m; or"
$mergeRange.Find.Execute($m4, $true, $false, $false, $false, $false, $true, 1, $false, $m4, 2) | Out-Null

$m5 = "(a + b) shall not exceed 4.5m for Residential developments, 5.0m for Commercial developments or 6.0m for Industrial or Warehouse developments."
$mergeRange.Find.Execute($m5, $true, $false, $false, $false, $false, $true, 1, $false, $m5, 2) | Out-Null

$m6 = "The structural beams shall also be included in the height computation."
$mergeRange.Find.Execute($m6, $true, $false, $false, $false, $false, $true, 1, $false, $m6, 2) | Out-Null

# ---------------------------------------------------------------
# 2. Replace the 3 inline pictures with hyperlinks to their source URLs
# ---------------------------------------------------------------

$urls = @(
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/GFA/GFA-27-Cable-chamber_integrated-substation_final.jpg?h=566&w=800",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/GFA/GFA-26A-Cable-chamber_standalone-substation_final.jpg?h=576&w=1000",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/GFA/GFA-26B-Cable-chamber_standalone-substation_final.jpg"
)

$n = $d.InlineShapes.Count
for ($i = 1; $i -le $n; $i++) {
    $shp = $d.InlineShapes.Item(1)
    $start = $shp.Range.Start
    $shp.Delete()
    $r = $d.Range($start, $start)
    $url = $urls[$i - 1]
    $d.Hyperlinks.Add($r, $url, "", "", $url) | Out-Null
}

Write-Output "done"
